# Auto-generated edit script: updates Kujata_Profits market-data cells per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2300
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 2160
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 2160
$ws.Range("M40").Value = -2825
$ws.Range("N40").Value = -2510
$ws.Range("H75").Value = 15314
$ws.Range("J75").Value = 15314
$ws.Range("L75").Value = 15314
$ws.Range("N75").Value = -17186
$ws.Range("H76").Value = 62503736
$ws.Range("J76").Value = 125003864
$ws.Range("L76").Value = 125003864
$ws.Range("N76").Value = -125004494
$ws.Range("H78").Value = 15314
$ws.Range("J78").Value = 15314
$ws.Range("L78").Value = 45942
$ws.Range("N78").Value = -55302
$ws.Range("H79").Value = 62503736
$ws.Range("J79").Value = 125003864
$ws.Range("L79").Value = 125003864
$ws.Range("N79").Value = -125006048
$ws.Range("H80").Value = 382.42105
$ws.Range("I80").Value = 205
$ws.Range("J80").Value = 579.55554
$ws.Range("K80").Value = 615
$ws.Range("L80").Value = 1738.66662
$ws.Range("M80").Value = 383
$ws.Range("N80").Value = -3734.66662
$ws.Range("H83").Value = 382.42105
$ws.Range("I83").Value = 205
$ws.Range("J83").Value = 579.55554
$ws.Range("K83").Value = 1845
$ws.Range("L83").Value = 5215.99986
$ws.Range("M83").Value = 3147
$ws.Range("N83").Value = -15199.99986
$ws.Range("H138").Value = 1232.1573
$ws.Range("J138").Value = 1726.5366
$ws.Range("L138").Value = 5179.6098
$ws.Range("N138").Value = -15459.6098

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4181.7803
$ws.Range("I32").Value = 3774.7368
$ws.Range("J32").Value = 9337.666999999999
$ws.Range("K32").Value = 3774.7368
$ws.Range("L32").Value = 9337.666999999999
$ws.Range("M32").Value = -3487.7368
$ws.Range("N32").Value = -9911.666999999999
$ws.Range("H96").Value = 17399.2
$ws.Range("J96").Value = 17399.2
$ws.Range("L96").Value = 17399.2
$ws.Range("N96").Value = -22891.2
$ws.Range("H107").Value = 28076
$ws.Range("J107").Value = 28076
$ws.Range("L107").Value = 28076
$ws.Range("N107").Value = -35756
$ws.Range("H109").Value = 41000
$ws.Range("J109").Value = 41000
$ws.Range("L109").Value = 41000
$ws.Range("N109").Value = -43774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3588.7778
$ws.Range("J20").Value = 3499.75
$ws.Range("L20").Value = 3499.75
$ws.Range("N20").Value = -3993.75
$ws.Range("H99").Value = 33334608
$ws.Range("I99").Value = 38462740
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 38462740
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = -38461242
$ws.Range("N99").Value = -4746
$ws.Range("H107").Value = 1571.4286
$ws.Range("I107").Value = 1340
$ws.Range("J107").Value = 2150
$ws.Range("K107").Value = 1340
$ws.Range("L107").Value = 2150
$ws.Range("M107").Value = 580
$ws.Range("N107").Value = -5990
$ws.Range("H123").Value = 50780
$ws.Range("J123").Value = 50780
$ws.Range("L123").Value = 50780
$ws.Range("N123").Value = -60580
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1753.6
$ws.Range("I31").Value = 1198.8889
$ws.Range("J31").Value = 2585.6667
$ws.Range("K31").Value = 1198.8889
$ws.Range("L31").Value = 2585.6667
$ws.Range("M31").Value = -903.8888999999999
$ws.Range("N31").Value = -3175.6667
$ws.Range("H34").Value = 1753.6
$ws.Range("I34").Value = 1198.8889
$ws.Range("J34").Value = 2585.6667
$ws.Range("K34").Value = 1198.8889
$ws.Range("L34").Value = 2585.6667
$ws.Range("M34").Value = -996.8888999999999
$ws.Range("N34").Value = -2989.6667
$ws.Range("H88").Value = 10000
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 10000
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H97").Value = 24900
$ws.Range("J97").Value = 24900
$ws.Range("L97").Value = 24900
$ws.Range("N97").Value = -26882
$ws.Range("H99").Value = 1839.0834
$ws.Range("I99").Value = 1824.4546
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1824.4546
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -326.4546
$ws.Range("N99").Value = -4996
$ws.Range("H122").Value = 893.8
$ws.Range("I122").Value = 893.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2681.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -231.3999999999996
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 1839.0834
$ws.Range("I126").Value = 1824.4546
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5473.3638
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3003.3638
$ws.Range("N126").Value = -10940
$ws.Range("H130").Value = 39600
$ws.Range("J130").Value = 39600
$ws.Range("L130").Value = 39600
$ws.Range("N130").Value = -49640
$ws.Range("H132").Value = 17260.285
$ws.Range("I132").Value = 35041.332
$ws.Range("J132").Value = 3924.5
$ws.Range("K132").Value = 105123.996
$ws.Range("L132").Value = 11773.5
$ws.Range("M132").Value = -102593.996
$ws.Range("N132").Value = -16833.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 477.21738
$ws.Range("I107").Value = 618.5
$ws.Range("J107").Value = 257.44446
$ws.Range("K107").Value = 618.5
$ws.Range("L107").Value = 257.44446
$ws.Range("M107").Value = 1301.5
$ws.Range("N107").Value = -4097.44446
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H133").Value = 42875.11
$ws.Range("J133").Value = 42875.11
$ws.Range("L133").Value = 42875.11
$ws.Range("N133").Value = -52995.11

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1508.3
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 1609.2222
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 1609.2222
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -2199.2222
$ws.Range("H27").Value = 1508.3
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 1609.2222
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 1609.2222
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -1823.2222
$ws.Range("H136").Value = 5062.207
$ws.Range("I136").Value = 6766.1665
$ws.Range("K136").Value = 20298.4995
$ws.Range("M136").Value = -17748.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 12726.333
$ws.Range("J52").Value = 12726.333
$ws.Range("L52").Value = 12726.333
$ws.Range("N52").Value = -13178.333
$ws.Range("H126").Value = 52632496
$ws.Range("I126").Value = 62500844
$ws.Range("K126").Value = 187502532
$ws.Range("M126").Value = -187500062
